$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.317.53'
$ws.Range("E2").Value = '  -0.13%  '

# Row 3
$ws.Range("D3").Value = '2.318.91'
$ws.Range("E3").Value = '  -2.13%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.60%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.10%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.628'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.26%  '

# Row 8
$ws.Range("E8").Value = '  +0.15%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.604'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.31%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.03'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.20%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0917'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.39%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.36'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.11%  '

# Row 13
$ws.Range("E13").Value = '  +0.63%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.977'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.48%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.59'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.09%  '

# Row 16
$ws.Range("D16").Value = '2.674.86'
$ws.Range("E16").Value = '  -1.87%  '

# Row 17
$ws.Range("D17").Value = '2.352.00'
$ws.Range("E17").Value = '  -0.45%  '

# Row 18
$ws.Range("D18").Value = '42.271.36'
$ws.Range("E18").Value = '  -0.37%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.04%  '

# Row 20
$ws.Range("E20").Value = '  -1.23%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '75.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.63%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.84%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '265.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.76%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.61%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.48%  '

# Row 26
$ws.Range("E26").Value = '  +0.37%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.19%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.57%  '

# Row 29
$ws.Range("E29").Value = '  +1.64%  '

# Row 30
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.81'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.80%  '

# Row 31
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '165.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.69%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0894'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.72%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.90'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.31%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.50%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.119'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +13.11%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.129'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.04%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.59'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.50%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0353'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.84%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.73'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.10%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.64'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -9.94%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '100.16'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +12.13%  '

# Row 42
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.46'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.49%  '

# Row 43
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.53'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.17%  '

# Row 44
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.232'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.60%  '

# Row 45
$ws.Range("E45").Value = '  +0.26%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.17'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.29%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '111.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.59%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.33%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.04'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.60%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.98'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.75%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.27'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.28%  '
